$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 <= old row 28
$ws.Cells.Item(27, 2).Value = 6883446
$ws.Cells.Item(27, 5).Value = 'KFC Komarno'
$ws.Cells.Item(27, 6).Value = 'FC Tatran Presov'
$ws.Cells.Item(27, 7).Value = 3
$ws.Cells.Item(27, 8).Value = 1
$ws.Cells.Item(27, 9).Value = 'H'
$ws.Cells.Item(27, 10).Value = 2.3
$ws.Cells.Item(27, 11).Value = 3.2
$ws.Cells.Item(27, 12).Value = 2.7
$ws.Cells.Item(27, 13).Value = 2.25
$ws.Cells.Item(27, 14).Value = 3.4
$ws.Cells.Item(27, 15).Value = 3.1
$ws.Cells.Item(27, 16).Value = -0.25
$ws.Cells.Item(27, 17).Value = 1.95
$ws.Cells.Item(27, 18).Value = 1.85
$ws.Cells.Item(27, 19).Value = 2.25
$ws.Cells.Item(27, 20).Value = 1.825
$ws.Cells.Item(27, 21).Value = 1.975
$ws.Cells.Item(27, 22).Value = 1.25
$ws.Cells.Item(27, 23).Value = -1
$ws.Cells.Item(27, 24).Value = -1
$ws.Cells.Item(27, 25).Value = 0.95
$ws.Cells.Item(27, 26).Value = -1
$ws.Cells.Item(27, 27).Value = 0.825
$ws.Cells.Item(27, 28).Value = -1

# Row 28 <= old row 29
$ws.Cells.Item(28, 2).Value = 6883849
$ws.Cells.Item(28, 5).Value = 'Puchov'
$ws.Cells.Item(28, 6).Value = 'OFK Malzenice'
$ws.Cells.Item(28, 7).Value = 2
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 9).Value = 'H'
$ws.Cells.Item(28, 10).Value = 1.4
$ws.Cells.Item(28, 11).Value = 4.333
$ws.Cells.Item(28, 12).Value = 6
$ws.Cells.Item(28, 13).Value = 1.5
$ws.Cells.Item(28, 14).Value = 4.5
$ws.Cells.Item(28, 15).Value = 6.5
$ws.Cells.Item(28, 16).Value = -1.25
$ws.Cells.Item(28, 17).Value = 1.95
$ws.Cells.Item(28, 18).Value = 1.85
$ws.Cells.Item(28, 19).Value = 2.75
$ws.Cells.Item(28, 20).Value = 1.75
$ws.Cells.Item(28, 21).Value = 1.95
$ws.Cells.Item(28, 22).Value = 0.5
$ws.Cells.Item(28, 23).Value = -1
$ws.Cells.Item(28, 24).Value = -1
$ws.Cells.Item(28, 25).Value = 0.95
$ws.Cells.Item(28, 26).Value = -1
$ws.Cells.Item(28, 27).Value = -1
$ws.Cells.Item(28, 28).Value = 0.95

# Row 29 <= old row 27
$ws.Cells.Item(29, 2).Value = 6884050
$ws.Cells.Item(29, 5).Value = 'MSK Povazska Bystrica'
$ws.Cells.Item(29, 6).Value = 'Slavoj Trebisov'
$ws.Cells.Item(29, 7).Value = 5
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).Value = 'H'
$ws.Cells.Item(29, 10).Value = 1.666
$ws.Cells.Item(29, 11).Value = 3.75
$ws.Cells.Item(29, 12).Value = 4
$ws.Cells.Item(29, 13).Value = 1.5
$ws.Cells.Item(29, 14).Value = 4.333
$ws.Cells.Item(29, 15).Value = 6.5
$ws.Cells.Item(29, 16).Value = -1.25
$ws.Cells.Item(29, 17).Value = 1.9
$ws.Cells.Item(29, 18).Value = 1.9
$ws.Cells.Item(29, 19).Value = 3
$ws.Cells.Item(29, 20).Value = 1.9
$ws.Cells.Item(29, 21).Value = 1.9
$ws.Cells.Item(29, 22).Value = 0.5
$ws.Cells.Item(29, 23).Value = -1
$ws.Cells.Item(29, 24).Value = -1
$ws.Cells.Item(29, 25).Value = 0.8999999999999999
$ws.Cells.Item(29, 26).Value = -1
$ws.Cells.Item(29, 27).Value = 0.8999999999999999
$ws.Cells.Item(29, 28).Value = -1

# Row 87 <= old row 89
$ws.Cells.Item(87, 2).Value = 6884106
$ws.Cells.Item(87, 5).Value = 'MSK Povazska Bystrica'
$ws.Cells.Item(87, 6).Value = 'MSK Zilina B'
$ws.Cells.Item(87, 7).Value = 2
$ws.Cells.Item(87, 8).Value = 1
$ws.Cells.Item(87, 9).Value = 'H'
$ws.Cells.Item(87, 10).Value = 1.8
$ws.Cells.Item(87, 11).Value = 3.7
$ws.Cells.Item(87, 12).Value = 3.5
$ws.Cells.Item(87, 13).Value = 2.5
$ws.Cells.Item(87, 14).Value = 3.6
$ws.Cells.Item(87, 15).Value = 2.6
$ws.Cells.Item(87, 16).Value = 0
$ws.Cells.Item(87, 17).Value = 1.85
$ws.Cells.Item(87, 18).Value = 1.95
$ws.Cells.Item(87, 19).Value = 3
$ws.Cells.Item(87, 20).Value = 1.9
$ws.Cells.Item(87, 21).Value = 1.9
$ws.Cells.Item(87, 22).Value = 1.5
$ws.Cells.Item(87, 23).Value = -1
$ws.Cells.Item(87, 24).Value = -1
$ws.Cells.Item(87, 25).Value = 0.8500000000000001
$ws.Cells.Item(87, 26).Value = -1
$ws.Cells.Item(87, 27).Value = 0
$ws.Cells.Item(87, 28).Value = 0

# Row 89 <= old row 87
$ws.Cells.Item(89, 2).Value = 6883864
$ws.Cells.Item(89, 5).Value = 'Puchov'
$ws.Cells.Item(89, 6).Value = 'KFC Komarno'
$ws.Cells.Item(89, 7).Value = 0
$ws.Cells.Item(89, 8).Value = 1
$ws.Cells.Item(89, 9).Value = 'A'
$ws.Cells.Item(89, 10).Value = 2.25
$ws.Cells.Item(89, 11).Value = 3.4
$ws.Cells.Item(89, 12).Value = 2.625
$ws.Cells.Item(89, 13).Value = 2.4
$ws.Cells.Item(89, 14).Value = 3.4
$ws.Cells.Item(89, 15).Value = 2.875
$ws.Cells.Item(89, 16).Value = -0.25
$ws.Cells.Item(89, 17).Value = 2.05
$ws.Cells.Item(89, 18).Value = 1.75
$ws.Cells.Item(89, 19).Value = 2.75
$ws.Cells.Item(89, 20).Value = 1.975
$ws.Cells.Item(89, 21).Value = 1.825
$ws.Cells.Item(89, 22).Value = -1
$ws.Cells.Item(89, 23).Value = -1
$ws.Cells.Item(89, 24).Value = 1.875
$ws.Cells.Item(89, 25).Value = -1
$ws.Cells.Item(89, 26).Value = 0.75
$ws.Cells.Item(89, 27).Value = -1
$ws.Cells.Item(89, 28).Value = 0.825

# Row 103 <= old row 105
$ws.Cells.Item(103, 2).Value = 6884124
$ws.Cells.Item(103, 5).Value = 'OFK Malzenice'
$ws.Cells.Item(103, 6).Value = 'Spisska Nova Ves'
$ws.Cells.Item(103, 7).Value = 3
$ws.Cells.Item(103, 8).Value = 1
$ws.Cells.Item(103, 9).Value = 'H'
$ws.Cells.Item(103, 10).Value = 2
$ws.Cells.Item(103, 11).Value = 3.4
$ws.Cells.Item(103, 12).Value = 3.1
$ws.Cells.Item(103, 13).Value = 2.05
$ws.Cells.Item(103, 14).Value = 3.6
$ws.Cells.Item(103, 15).Value = 3.5
$ws.Cells.Item(103, 16).Value = -0.5
$ws.Cells.Item(103, 17).Value = 2
$ws.Cells.Item(103, 18).Value = 1.8
$ws.Cells.Item(103, 19).Value = 2.25
$ws.Cells.Item(103, 20).Value = 2
$ws.Cells.Item(103, 21).Value = 1.8
$ws.Cells.Item(103, 22).Value = 1.05
$ws.Cells.Item(103, 23).Value = -1
$ws.Cells.Item(103, 24).Value = -1
$ws.Cells.Item(103, 25).Value = 1
$ws.Cells.Item(103, 26).Value = -1
$ws.Cells.Item(103, 27).Value = 1
$ws.Cells.Item(103, 28).Value = -1

# Row 104 <= old row 103
$ws.Cells.Item(104, 2).Value = 6884112
$ws.Cells.Item(104, 5).Value = 'STK Samorin'
$ws.Cells.Item(104, 6).Value = 'Spartak Myjava'
$ws.Cells.Item(104, 7).Value = 2
$ws.Cells.Item(104, 8).Value = 2
$ws.Cells.Item(104, 9).Value = 'D'
$ws.Cells.Item(104, 10).Value = 2.625
$ws.Cells.Item(104, 11).Value = 3.4
$ws.Cells.Item(104, 12).Value = 2.4
$ws.Cells.Item(104, 13).Value = 3.1
$ws.Cells.Item(104, 14).Value = 3.75
$ws.Cells.Item(104, 15).Value = 2.1
$ws.Cells.Item(104, 16).Value = 0.25
$ws.Cells.Item(104, 17).Value = 1.95
$ws.Cells.Item(104, 18).Value = 1.85
$ws.Cells.Item(104, 19).Value = 2.5
$ws.Cells.Item(104, 20).Value = 1.775
$ws.Cells.Item(104, 21).Value = 1.925
$ws.Cells.Item(104, 22).Value = -1
$ws.Cells.Item(104, 23).Value = 2.75
$ws.Cells.Item(104, 24).Value = -1
$ws.Cells.Item(104, 25).Value = 0.475
$ws.Cells.Item(104, 26).Value = -0.5
$ws.Cells.Item(104, 27).Value = 0.7749999999999999
$ws.Cells.Item(104, 28).Value = -1

# Row 105 <= old row 104
$ws.Cells.Item(105, 2).Value = 6878191
$ws.Cells.Item(105, 5).Value = 'FC Petrzalka'
$ws.Cells.Item(105, 6).Value = 'KFC Komarno'
$ws.Cells.Item(105, 7).Value = 3
$ws.Cells.Item(105, 8).Value = 1
$ws.Cells.Item(105, 9).Value = 'H'
$ws.Cells.Item(105, 10).Value = 3
$ws.Cells.Item(105, 11).Value = 3.6
$ws.Cells.Item(105, 12).Value = 2.05
$ws.Cells.Item(105, 13).Value = 2.8
$ws.Cells.Item(105, 14).Value = 3.75
$ws.Cells.Item(105, 15).Value = 2.25
$ws.Cells.Item(105, 16).Value = 0.25
$ws.Cells.Item(105, 17).Value = 1.825
$ws.Cells.Item(105, 18).Value = 1.975
$ws.Cells.Item(105, 19).Value = 2.5
$ws.Cells.Item(105, 20).Value = 1.925
$ws.Cells.Item(105, 21).Value = 1.875
$ws.Cells.Item(105, 22).Value = 1.8
$ws.Cells.Item(105, 23).Value = -1
$ws.Cells.Item(105, 24).Value = -1
$ws.Cells.Item(105, 25).Value = 0.825
$ws.Cells.Item(105, 26).Value = -1
$ws.Cells.Item(105, 27).Value = 0.925
$ws.Cells.Item(105, 28).Value = -1

# Row 141 <= old row 142
$ws.Cells.Item(141, 2).Value = 6884078
$ws.Cells.Item(141, 5).Value = 'Spisska Nova Ves'
$ws.Cells.Item(141, 6).Value = 'FK Pohronie'
$ws.Cells.Item(141, 7).Value = 0
$ws.Cells.Item(141, 8).Value = 1
$ws.Cells.Item(141, 9).Value = 'A'
$ws.Cells.Item(141, 10).Value = 4.333
$ws.Cells.Item(141, 11).Value = 3.75
$ws.Cells.Item(141, 12).Value = 1.615
$ws.Cells.Item(141, 13).Value = 5
$ws.Cells.Item(141, 14).Value = 4.2
$ws.Cells.Item(141, 15).Value = 1.6
$ws.Cells.Item(141, 16).Value = 0.75
$ws.Cells.Item(141, 17).Value = 1.95
$ws.Cells.Item(141, 18).Value = 1.75
$ws.Cells.Item(141, 19).Value = 2.75
$ws.Cells.Item(141, 20).Value = 1.875
$ws.Cells.Item(141, 21).Value = 1.925
$ws.Cells.Item(141, 22).Value = -1
$ws.Cells.Item(141, 23).Value = -1
$ws.Cells.Item(141, 24).Value = 0.6000000000000001
$ws.Cells.Item(141, 25).Value = -0.5
$ws.Cells.Item(141, 26).Value = 0.375
$ws.Cells.Item(141, 27).Value = -1
$ws.Cells.Item(141, 28).Value = 0.925

# Row 142 <= old row 141
$ws.Cells.Item(142, 2).Value = 6884076
$ws.Cells.Item(142, 5).Value = 'Slavoj Trebisov'
$ws.Cells.Item(142, 6).Value = 'MSK Povazska Bystrica'
$ws.Cells.Item(142, 7).Value = 1
$ws.Cells.Item(142, 8).Value = 2
$ws.Cells.Item(142, 9).Value = 'A'
$ws.Cells.Item(142, 10).Value = 2.625
$ws.Cells.Item(142, 11).Value = 3.4
$ws.Cells.Item(142, 12).Value = 2.3
$ws.Cells.Item(142, 13).Value = 2.625
$ws.Cells.Item(142, 14).Value = 3.4
$ws.Cells.Item(142, 15).Value = 2.625
$ws.Cells.Item(142, 16).Value = 0
$ws.Cells.Item(142, 17).Value = 1.9
$ws.Cells.Item(142, 18).Value = 1.9
$ws.Cells.Item(142, 19).Value = 2.25
$ws.Cells.Item(142, 20).Value = 1.8
$ws.Cells.Item(142, 21).Value = 2
$ws.Cells.Item(142, 22).Value = -1
$ws.Cells.Item(142, 23).Value = -1
$ws.Cells.Item(142, 24).Value = 1.625
$ws.Cells.Item(142, 25).Value = -1
$ws.Cells.Item(142, 26).Value = 0.8999999999999999
$ws.Cells.Item(142, 27).Value = 0.8
$ws.Cells.Item(142, 28).Value = -1

# Row 152 <= old row 153
$ws.Cells.Item(152, 2).Value = 6884162
$ws.Cells.Item(152, 5).Value = 'STK Samorin'
$ws.Cells.Item(152, 6).Value = 'Dolny Kubin'
$ws.Cells.Item(152, 7).Value = 4
$ws.Cells.Item(152, 8).Value = 0
$ws.Cells.Item(152, 9).Value = 'H'
$ws.Cells.Item(152, 10).Value = 1.571
$ws.Cells.Item(152, 11).Value = 4
$ws.Cells.Item(152, 12).Value = 5
$ws.Cells.Item(152, 13).Value = 1.571
$ws.Cells.Item(152, 14).Value = 4.5
$ws.Cells.Item(152, 15).Value = 5.25
$ws.Cells.Item(152, 16).Value = -1
$ws.Cells.Item(152, 17).Value = 1.8
$ws.Cells.Item(152, 18).Value = 2
$ws.Cells.Item(152, 19).Value = 3
$ws.Cells.Item(152, 20).Value = 1.8
$ws.Cells.Item(152, 21).Value = 2
$ws.Cells.Item(152, 22).Value = 0.571
$ws.Cells.Item(152, 23).Value = -1
$ws.Cells.Item(152, 24).Value = -1
$ws.Cells.Item(152, 25).Value = 0.8
$ws.Cells.Item(152, 26).Value = -1
$ws.Cells.Item(152, 27).Value = 0.8
$ws.Cells.Item(152, 28).Value = -1

# Row 153 <= old row 152
$ws.Cells.Item(153, 2).Value = 6883468
$ws.Cells.Item(153, 5).Value = 'MSK Zilina B'
$ws.Cells.Item(153, 6).Value = 'FC Tatran Presov'
$ws.Cells.Item(153, 7).Value = 1
$ws.Cells.Item(153, 8).Value = 2
$ws.Cells.Item(153, 9).Value = 'A'
$ws.Cells.Item(153, 10).Value = 4.8
$ws.Cells.Item(153, 11).Value = 3.8
$ws.Cells.Item(153, 12).Value = 1.6
$ws.Cells.Item(153, 13).Value = 5.25
$ws.Cells.Item(153, 14).Value = 4.2
$ws.Cells.Item(153, 15).Value = 1.615
$ws.Cells.Item(153, 16).Value = 1
$ws.Cells.Item(153, 17).Value = 1.8
$ws.Cells.Item(153, 18).Value = 2
$ws.Cells.Item(153, 19).Value = 3.25
$ws.Cells.Item(153, 20).Value = 1.975
$ws.Cells.Item(153, 21).Value = 1.825
$ws.Cells.Item(153, 22).Value = -1
$ws.Cells.Item(153, 23).Value = -1
$ws.Cells.Item(153, 24).Value = 0.615
$ws.Cells.Item(153, 25).Value = 0
$ws.Cells.Item(153, 26).Value = 0
$ws.Cells.Item(153, 27).Value = -0.5
$ws.Cells.Item(153, 28).Value = 0.4125

# Row 154 <= old row 155
$ws.Cells.Item(154, 2).Value = 6883469
$ws.Cells.Item(154, 5).Value = 'FK Humenne'
$ws.Cells.Item(154, 6).Value = 'STK Samorin'
$ws.Cells.Item(154, 7).Value = 1
$ws.Cells.Item(154, 8).Value = 0
$ws.Cells.Item(154, 9).Value = 'H'
$ws.Cells.Item(154, 10).Value = 1.727
$ws.Cells.Item(154, 11).Value = 3.75
$ws.Cells.Item(154, 12).Value = 3.75
$ws.Cells.Item(154, 13).Value = 1.533
$ws.Cells.Item(154, 14).Value = 4.5
$ws.Cells.Item(154, 15).Value = 5.5
$ws.Cells.Item(154, 16).Value = -1
$ws.Cells.Item(154, 17).Value = 1.875
$ws.Cells.Item(154, 18).Value = 1.925
$ws.Cells.Item(154, 19).Value = 2.75
$ws.Cells.Item(154, 20).Value = 1.85
$ws.Cells.Item(154, 21).Value = 1.95
$ws.Cells.Item(154, 22).Value = 0.5329999999999999
$ws.Cells.Item(154, 23).Value = -1
$ws.Cells.Item(154, 24).Value = -1
$ws.Cells.Item(154, 25).Value = 0
$ws.Cells.Item(154, 26).Value = 0
$ws.Cells.Item(154, 27).Value = -1
$ws.Cells.Item(154, 28).Value = 0.95

# Row 155 <= old row 156
$ws.Cells.Item(155, 2).Value = 6883470
$ws.Cells.Item(155, 5).Value = 'FC Tatran Presov'
$ws.Cells.Item(155, 6).Value = 'FK Pohronie'
$ws.Cells.Item(155, 7).Value = 2
$ws.Cells.Item(155, 8).Value = 0
$ws.Cells.Item(155, 9).Value = 'H'
$ws.Cells.Item(155, 10).Value = 1.333
$ws.Cells.Item(155, 11).Value = 5
$ws.Cells.Item(155, 12).Value = 7
$ws.Cells.Item(155, 13).Value = 1.3
$ws.Cells.Item(155, 14).Value = 6
$ws.Cells.Item(155, 15).Value = 9
$ws.Cells.Item(155, 16).Value = -1.75
$ws.Cells.Item(155, 17).Value = 1.975
$ws.Cells.Item(155, 18).Value = 1.825
$ws.Cells.Item(155, 19).Value = 2.75
$ws.Cells.Item(155, 20).Value = 1.875
$ws.Cells.Item(155, 21).Value = 1.925
$ws.Cells.Item(155, 22).Value = 0.3
$ws.Cells.Item(155, 23).Value = -1
$ws.Cells.Item(155, 24).Value = -1
$ws.Cells.Item(155, 25).Value = 0.4875
$ws.Cells.Item(155, 26).Value = -0.5
$ws.Cells.Item(155, 27).Value = -1
$ws.Cells.Item(155, 28).Value = 0.925

# Row 156 <= old row 154
$ws.Cells.Item(156, 2).Value = 6884081
$ws.Cells.Item(156, 5).Value = 'Spisska Nova Ves'
$ws.Cells.Item(156, 6).Value = 'Spartak Myjava'
$ws.Cells.Item(156, 7).Value = 0
$ws.Cells.Item(156, 8).Value = 0
$ws.Cells.Item(156, 9).Value = 'D'
$ws.Cells.Item(156, 10).Value = 4.333
$ws.Cells.Item(156, 11).Value = 3.75
$ws.Cells.Item(156, 12).Value = 1.615
$ws.Cells.Item(156, 13).Value = 5.25
$ws.Cells.Item(156, 14).Value = 3.75
$ws.Cells.Item(156, 15).Value = 1.666
$ws.Cells.Item(156, 16).Value = 0.75
$ws.Cells.Item(156, 17).Value = 1.875
$ws.Cells.Item(156, 18).Value = 1.925
$ws.Cells.Item(156, 19).Value = 2.5
$ws.Cells.Item(156, 20).Value = 1.9
$ws.Cells.Item(156, 21).Value = 1.9
$ws.Cells.Item(156, 22).Value = -1
$ws.Cells.Item(156, 23).Value = 2.75
$ws.Cells.Item(156, 24).Value = -1
$ws.Cells.Item(156, 25).Value = 0.875
$ws.Cells.Item(156, 26).Value = -1
$ws.Cells.Item(156, 27).Value = -1
$ws.Cells.Item(156, 28).Value = 0.8999999999999999

